# Commit: "remain only one choice of gas (fossil fuel)"
#
# The "supply" sheet listed three gas supply technologies
# (gas_standard, gas_individuell, gas_erneuerbar). Only one choice of gas
# (the fossil-fuel "gas_standard") should remain, so the two extra gas
# rows ("gas_individuell" and "gas_erneuerbar", originally rows 6 and 7)
# are removed. Every row below shifts up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supply")

# Preserve the text of the cell-note comments that live on the rows which
# will shift up, so they can be re-anchored on their new cells after the
# row delete (the engine moves cell data/formulas but not legacy notes).
$textE8  = $ws.Range("E8").Comment.Text()
$textG8  = $ws.Range("G8").Comment.Text()
$textE9  = $ws.Range("E9").Comment.Text()
$textG10 = $ws.Range("G10").Comment.Text()

# Drop the now-stale comments before the shift so they don't linger on
# the wrong cells once the rows move.
$ws.Range("E8").Comment.Delete()
$ws.Range("G8").Comment.Delete()
$ws.Range("E9").Comment.Delete()
$ws.Range("G10").Comment.Delete()

# Remove the two "gas_individuell" / "gas_erneuerbar" rows; everything
# beneath (pallet, oil, district_heating) moves up two rows.
$ws.Rows.Item(6).Resize(2).Delete()

# Re-create the comments on their new, shifted-up cells.
$ws.Range("E6").AddComment($textE8)
$ws.Range("G6").AddComment($textG8)
$ws.Range("E7").AddComment($textE9)
$ws.Range("G8").AddComment($textG10)

$ws.Activate()
$ws.Range("G5").Select()
